# Rename the failed-status label from "Still Incorrect" to "Didn't get it"
# inside the quoted list: "Strike!","Spare!","Still Incorrect"

$d = $word.ActiveDocument

# Locate "Still Incorrect" so we know exactly where the replacement text
# begins (this is also where the _GoBack bookmark should end up, mirroring
# where Word leaves it after an in-place retype).
$rng = $d.Content
$rng.Find.Execute("Still Incorrect", $false, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0) | Out-Null
$splitPos = $rng.Start

# Replace the found text with the new wording. Assigning .Text (rather than
# going through Find/Replace) keeps the literal straight apostrophe instead
# of letting autocorrect turn it into a curly one.
$rng.Text = "Didn't get it"

# Re-anchor the _GoBack bookmark at the point where the new text begins.
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
